# New crime data collected - weekly CompStat update (67th Precinct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/report dates) ---
$ws.Range("A8").Value = "Volume 32   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# --- Row 14: Murder ---
$ws.Range("F14").Value = 1
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("I14").Value = 6
$ws.Range("K14").Value = -25
$ws.Range("L14").Value = -25
$ws.Range("M14").Value = -33.333333333333
$ws.Range("N14").Value = -70

# --- Row 15: Rape ---
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 19
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = 5.555555555555
$ws.Range("L15").Value = -9.523809523809
$ws.Range("M15").Value = 5.555555555555
$ws.Range("N15").Value = -59.574468085106

# --- Row 16: Robbery ---
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 14.285714285714
$ws.Range("I16").Value = 138
$ws.Range("J16").Value = 143
$ws.Range("K16").Value = -3.496503496503
$ws.Range("L16").Value = 4.545454545454
$ws.Range("M16").Value = -31.683168316831
$ws.Range("N16").Value = -88.364249578414

# --- Row 17: Fel. Assault ---
$ws.Range("C17").Value = 22
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = -4.347826086956
$ws.Range("F17").Value = 65
$ws.Range("G17").Value = 75
$ws.Range("H17").Value = -13.333333333333
$ws.Range("I17").Value = 359
$ws.Range("J17").Value = 403
$ws.Range("K17").Value = -10.918114143920
$ws.Range("L17").Value = 9.451219512195
$ws.Range("M17").Value = 63.926940639269
$ws.Range("N17").Value = -39.663865546218

# --- Row 18: Burglary ---
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 86
$ws.Range("J18").Value = 108
$ws.Range("K18").Value = -20.370370370370
$ws.Range("L18").Value = -17.307692307692
$ws.Range("M18").Value = -52.747252747252
$ws.Range("N18").Value = -91.518737672583

# --- Row 19: Gr. Larceny ---
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 22.222222222222
$ws.Range("F19").Value = 50
$ws.Range("H19").Value = 6.382978723404
$ws.Range("I19").Value = 325
$ws.Range("J19").Value = 294
$ws.Range("K19").Value = 10.544217687074
$ws.Range("L19").Value = -6.069364161849
$ws.Range("M19").Value = 22.641509433962
$ws.Range("N19").Value = -10.714285714285

# --- Row 20: G.L.A. ---
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = -71.428571428571
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = -21.739130434782
$ws.Range("I20").Value = 95
$ws.Range("J20").Value = 107
$ws.Range("K20").Value = -11.214953271028
$ws.Range("L20").Value = -24.603174603174
$ws.Range("M20").Value = -33.566433566433
$ws.Range("N20").Value = -90.960989533777

# --- Row 21: TOTAL ---
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = -8.510638297872
$ws.Range("F21").Value = 164
$ws.Range("G21").Value = 179
$ws.Range("H21").Value = -8.379888268156
$ws.Range("I21").Value = 1028
$ws.Range("J21").Value = 1081
$ws.Range("K21").Value = -4.902867715078
$ws.Range("L21").Value = -3.474178403755
$ws.Range("M21").Value = -0.963391136801
$ws.Range("N21").Value = -75.964461070844

# --- Row 22: Transit ---
$ws.Range("L22").Value = -60

# --- Row 24: Petit Larceny ---
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -23.809523809523
$ws.Range("F24").Value = 80
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = -25.233644859813
$ws.Range("I24").Value = 528
$ws.Range("J24").Value = 629
$ws.Range("K24").Value = -16.057233704292
$ws.Range("L24").Value = -12.582781456953
$ws.Range("M24").Value = 14.038876889848

# --- Row 25: Retail Theft ---
$ws.Range("C25").Value = 2
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 97
$ws.Range("J25").Value = 90
$ws.Range("K25").Value = 7.777777777777
$ws.Range("L25").Value = -17.094017094017

# --- Row 26: Misd. Assault ---
$ws.Range("C26").Value = 30
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = 66.666666666666
$ws.Range("F26").Value = 84
$ws.Range("G26").Value = 70
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 449
$ws.Range("J26").Value = 445
$ws.Range("K26").Value = 0.898876404494
$ws.Range("L26").Value = 19.098143236074
$ws.Range("M26").Value = -1.535087719298

# --- Row 27: UCR Rape* ---
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -75
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -71.428571428571
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 29
$ws.Range("K27").Value = -27.586206896551
$ws.Range("L27").Value = -22.222222222222

# --- Row 28: Other Sex Crimes ---
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 46
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 15
$ws.Range("L28").Value = 17.948717948717

# --- Row 29: Shooting Vic. ---
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = -33.333333333333
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = -37.5
$ws.Range("I29").Value = 22
$ws.Range("J29").Value = 20
$ws.Range("K29").Value = 10
$ws.Range("L29").Value = -4.347826086956
$ws.Range("M29").Value = -26.666666666666
$ws.Range("N29").Value = -77.083333333333

# --- Row 30: Shooting Inc. ---
$ws.Range("I30").Value = 17
$ws.Range("J30").Value = 17
$ws.Range("L30").Value = -10.526315789473
$ws.Range("M30").Value = -34.615384615384
$ws.Range("N30").Value = -81.318681318681

# --- Column H width auto-fit shrink (content no longer needs extra width) ---
# Target stored width is 6.168446 (matches the other data columns); the
# host quantizes column widths to an MDW-7 pixel grid, so 5.4 chars is the
# input that lands closest to that target after round-tripping.
$ws.Columns("H").ColumnWidth = 5.4
